$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (not auto-converted numbers) for D/E columns, matching source data
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '40.140.77'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '2.223.91'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '294.33'
$ws.Range("E5").Value = '  +1.48%  '

$ws.Range("D6").Value = '87.91'
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").Value = '0.514'
$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").Value = '30.74'
$ws.Range("E10").Value = '  +0.07%  '

$ws.Range("D11").Value = '50.88'
$ws.Range("E11").Value = '  +6.45%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("E13").Value = '  +3.46%  '

$ws.Range("E14").Value = '  -0.53%  '

$ws.Range("D15").Value = '2.549.06'
$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("D16").Value = '13.85'
$ws.Range("E16").Value = '  -1.21%  '

$ws.Range("D17").Value = '2.237.22'
$ws.Range("E17").Value = '  +1.70%  '

$ws.Range("D18").Value = '0.738'
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("D19").Value = '40.050.63'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  -5.33%  '

$ws.Range("D22").Value = '5.78'
$ws.Range("E22").Value = '  -0.56%  '

$ws.Range("D23").Value = '65.72'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Value = '236.25'
$ws.Range("E24").Value = '  +0.65%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").Value = '2.49'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").Value = '23.22'
$ws.Range("E28").Value = '  +2.60%  '

$ws.Range("E29").Value = '  +1.18%  '

$ws.Range("E30").Value = '  -8.46%  '

$ws.Range("D31").Value = '158.89'
$ws.Range("E31").Value = '  +3.48%  '

$ws.Range("D32").Value = '31.86'
$ws.Range("E32").Value = '  -1.06%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").Value = '4.98'
$ws.Range("E34").Value = '  +0.32%  '

$ws.Range("D35").Value = '3.01'
$ws.Range("E35").Value = '  +6.45%  '

$ws.Range("E36").Value = '  -0.60%  '

$ws.Range("E37").Value = '  -2.88%  '

$ws.Range("D38").Value = '0.114'
$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("D39").Value = '1.76'
$ws.Range("E39").Value = '  +2.95%  '

$ws.Range("D40").Value = '0.0996'
$ws.Range("E40").Value = '  -0.66%  '

$ws.Range("D41").Value = '15.60'
$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("D42").Value = '2.083.35'

$ws.Range("D43").Value = '3.76'
$ws.Range("E43").Value = '  -2.38%  '

$ws.Range("D44").Value = '19.19'
$ws.Range("E44").Value = '  +7.66%  '

$ws.Range("E45").Value = '  +1.54%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  +2.32%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '1.89'
$ws.Range("E48").Value = '  -13.66%  '

$ws.Range("D49").Value = '2.450.70'
$ws.Range("E49").Value = '  +0.80%  '

$ws.Range("E50").Value = '  +2.08%  '

$ws.Range("E51").Value = '  +4.14%  '

# Restore default (unformatted) style now that values are stored as text
$textRange.ClearFormats()
